$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Starting" value (E2 label) in F2
$ws.Range("F2").Value = 1522529

# Update the "Test" value (E3 label) in F3
$ws.Range("F3").Value = 1338209

# Move the active selection to F3, matching the committed cursor position
$ws.Range("F3").Select()
